# The commit swaps the two theme parts of the deck: the theme that was
# carrying the "Office Theme" colour scheme (dk1/lt1/dk2/lt2/accent1-6/
# hlink/folHlink) and the theme carrying the "Integral" colour scheme trade
# places, so the presentation's live theme (the one driving SlideMaster1 /
# all 25 slides) ends up with the plain "Office" palette instead of
# "Integral".
#
# PowerPoint's COM model doesn't expose a "swap underlying theme part" verb,
# so the supported way to repaint a theme in this host is per-colour via
# ThemeColorScheme.Colors(i).RGB (fonts/format scheme are unaffected here -
# they were already identical between the two themes).
#
# Colour order for ThemeColorScheme.Colors(1..12):
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
#
# Target values are the old "Office Theme" srgbClr values (converted to the
# VBA-style 0x00BBGGRR integer RGB() encoding):
#   dk1      000000 -> 0
#   lt1      FFFFFF -> 16777215
#   dk2      44546A -> 6968388
#   lt2      E7E6E6 -> 15132391
#   accent1  5B9BD5 -> 13998939
#   accent2  ED7D31 -> 3243501
#   accent3  A5A5A5 -> 10855845
#   accent4  FFC000 -> 49407
#   accent5  4472C4 -> 12874308
#   accent6  70AD47 -> 4697456
#   hlink    0563C1 -> 12673797
#   folHlink 954F72 -> 7491477

$p = $ppt.ActivePresentation
$design = $p.Designs.Item(1)
$master = $design.SlideMaster
$theme = $master.Theme
$tcs = $theme.ThemeColorScheme

$targetRgb = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

for ($i = 1; $i -le 12; $i++) {
    $color = $tcs.Colors($i)
    $color.RGB = $targetRgb[$i - 1]
}
